# Minimal changes to base: update login credentials and add a hyperlink
# on the password cell (B2) matching the username cell's hyperlink style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the displayed credential text ---
$ws.Range("A2").Value = "matschie@testleaf.com"
$ws.Range("B2").Value = "SelBootCamp@123"

# --- Refresh the username hyperlink to point at the new address ---
$ws.Range("A2").Hyperlinks.Delete()
$hUser = $ws.Hyperlinks.Add($ws.Range("A2"), "mailto:matschie@testleaf.com")
$ws.Range("A2").Style = "Hyperlink"

# --- Add a new hyperlink on the password cell, matching username style ---
$hPass = $ws.Hyperlinks.Add($ws.Range("B2"), "mailto:matschie@testleaf.com")
$ws.Range("B2").Style = "Hyperlink"

Write-Output "done"
